# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.805.16"
$ws.Range("E2").Value = "  -3.37%  "
$ws.Range("D3").Value = "2.906.94"
$ws.Range("E3").Value = "  -4.13%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'586.98"
$ws.Range("E5").Value = "  -1.10%  "
$ws.Range("D6").Value = "'144.64"
$ws.Range("E6").Value = "  -5.81%  "
$ws.Range("E8").Value = "  -2.80%  "
$ws.Range("D9").Value = "2.907.76"
$ws.Range("E9").Value = "  -3.93%  "
$ws.Range("D10").Value = "'6.68"
$ws.Range("E10").Value = "  -1.78%  "
$ws.Range("E11").Value = "  -4.91%  "
$ws.Range("E12").Value = "  -3.98%  "
$ws.Range("D13").Value = "'0.0000226"
$ws.Range("E13").Value = "  -3.38%  "
$ws.Range("D14").Value = "'33.43"
$ws.Range("E14").Value = "  -6.50%  "
$ws.Range("E15").Value = "  +1.61%  "
$ws.Range("D16").Value = "3.390.13"
$ws.Range("E16").Value = "  -4.10%  "
$ws.Range("D17").Value = "60.773.54"
$ws.Range("E17").Value = "  -3.33%  "
$ws.Range("D18").Value = "'6.71"
$ws.Range("E18").Value = "  -5.14%  "
$ws.Range("D19").Value = "2.907.33"
$ws.Range("E19").Value = "  -4.15%  "
$ws.Range("D20").Value = "'427.87"
$ws.Range("E20").Value = "  -5.49%  "
$ws.Range("D21").Value = "'13.54"
$ws.Range("E21").Value = "  -5.00%  "
$ws.Range("D22").Value = "'0.680"
$ws.Range("E22").Value = "  -2.62%  "
$ws.Range("D23").Value = "'7.10"
$ws.Range("E23").Value = "  -5.46%  "
$ws.Range("D24").Value = "'80.51"
$ws.Range("E24").Value = "  -3.24%  "
$ws.Range("D25").Value = "'10.86"
$ws.Range("E25").Value = "  -3.15%  "
$ws.Range("E26").Value = "  -2.41%  "
$ws.Range("D27").Value = "'11.91"
$ws.Range("E27").Value = "  -4.26%  "
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("D30").Value = "'7.22"
$ws.Range("E30").Value = "  -3.74%  "
$ws.Range("E31").Value = "  -3.16%  "
$ws.Range("D32").Value = "'2.17"
$ws.Range("E32").Value = "  -3.83%  "
$ws.Range("E33").Value = "  -4.14%  "
$ws.Range("D34").Value = "'0.106"
$ws.Range("E34").Value = "  -3.70%  "
$ws.Range("D35").Value = "0.0₃0876"
$ws.Range("E35").Value = "  +2.03%  "
$ws.Range("E36").Value = "  -3.30%  "
$ws.Range("D37").Value = "'5.61"
$ws.Range("E37").Value = "  -5.46%  "
$ws.Range("D38").Value = "'3.02"
$ws.Range("E38").Value = "  -5.56%  "
$ws.Range("E39").Value = "  -2.78%  "
$ws.Range("E40").Value = "  -2.14%  "
$ws.Range("E41").Value = "  -4.73%  "
$ws.Range("D42").Value = "'8.61"
$ws.Range("E42").Value = "  -5.93%  "
$ws.Range("D43").Value = "'0.296"
$ws.Range("E43").Value = "  -3.74%  "
$ws.Range("D44").Value = "'41.44"
$ws.Range("E44").Value = "  -6.19%  "
$ws.Range("B45").Value = "Bittensor"
$ws.Range("C45").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D45").Value = "'378.10"
$ws.Range("E45").Value = "  -3.24%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "'0.0350"
$ws.Range("E46").Value = "  -2.90%  "
$ws.Range("D47").Value = "2.692.29"
$ws.Range("E47").Value = "  -1.10%  "
$ws.Range("D48").Value = "'132.41"
$ws.Range("E48").Value = "  -0.94%  "
$ws.Range("D50").Value = "'24.38"
$ws.Range("E50").Value = "  -1.91%  "
$ws.Range("E51").Value = "  -2.55%  "

# Reset style to Normal for text-forced numeric-looking cells to avoid quotePrefix style drift
foreach ($ref in @("D5","D6","D10","D13","D14","D18","D20","D21","D22","D23","D24","D25","D27","D30","D32","D34","D37","D38","D42","D43","D44","D45","D46","D48","D50")) {
    $ws.Range($ref).Style = "Normal"
}
